# AltController - Strings_for_translation.xlsx
# Add new translation rows for the "Snooze" feature (v1.94)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new rows just above "Situations_Title" (old row 380) for the
#    new Situations_SnoozeLabel / Situations_SnoozeToolTip strings.
# ---------------------------------------------------------------------------
$ws.Rows("380:381").Insert()

$ws.Range("A380").Value = "Situations_SnoozeLabel"
$ws.Range("C380").Value = "Edit situations"
$ws.Range("D380").Value = "Snooze profile when using this app"

$ws.Range("A381").Value = "Situations_SnoozeToolTip"
$ws.Range("C381").Value = "Edit situations"
$ws.Range("D381").Value = "Don't perform any actions while this app is active"

# ---------------------------------------------------------------------------
# 2) Insert one new row just above "String_Standard_pointer" (old row 501,
#    now shifted to row 503 because of the two rows inserted above) for the
#    new String_snooze string.
# ---------------------------------------------------------------------------
$ws.Rows("503:503").Insert()

$ws.Range("A503").Value = "String_snooze"
$ws.Range("C503").Value = "Strings"
$ws.Range("D503").Value = "snooze"

# ---------------------------------------------------------------------------
# Leave the selection where the last edit happened, similar to what Excel
# would record after typing these new values in.
# ---------------------------------------------------------------------------
$ws.Range("A503:D503").Select()
